$wb = $excel.ActiveWorkbook
$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet = $wb.Worksheets.Item("used")

# The name to move from the pool (Sheet1) to the used list (used sheet)
$name = $namesSheet.Range("A2").Value2
$fileName = "ChatGPT Image 2026年1月21日 16_28_28.png"
$usedAt = "2026-01-21 16:29:33"

# Remove the used name from the pool, shifting remaining names up
$namesSheet.Rows.Item(2).Delete()

# Find the next empty row in the used sheet and append the record
$lastRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$usedSheet.Cells.Item($newRow, 1).Value = $name
$usedSheet.Cells.Item($newRow, 2).Value = $fileName
$usedSheet.Cells.Item($newRow, 3).Value = $usedAt
